# "incluindo many to many"
# Turn the Vendas sheet into a proper many-to-many junction table:
# add an id_venda primary-key column in front of id_cliente / id_produto,
# and renumber/replace the old "quantidade" column's data accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vendas")

# New headers: A=id_venda, B=id_cliente, C=id_produto
$ws.Cells.Item(1, 1).Value = "id_venda"
$ws.Cells.Item(1, 2).Value = "id_cliente"
$ws.Cells.Item(1, 3).Value = "id_produto"

# New data rows (id_venda, id_cliente, id_produto)
$data = @(
    @(1, 1, 2),
    @(1, 1, 1),
    @(2, 2, 3),
    @(3, 3, 2),
    @(4, 5, 5),
    @(5, 4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Vendas becomes the active/selected sheet, with D11 selected.
$ws.Activate() | Out-Null
$ws.Range("D11").Select() | Out-Null
